$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(2, 4).Range.Text = "75"
$t.Cell(3, 4).Range.Text = "6 (8.0)"
$t.Cell(3, 6).Range.Text = "0.445"
$t.Cell(3, 8).Range.Text = "0.539"
$t.Cell(4, 4).Range.Text = "69 (92.0)"
$t.Cell(5, 4).Range.Text = "26 (34.7)"
$t.Cell(5, 6).Range.Text = "0.747"
$t.Cell(6, 4).Range.Text = "49 (65.3)"
$t.Cell(7, 4).Range.Text = "46 (61.3)"
$t.Cell(7, 6).Range.Text = "0.135"
$t.Cell(7, 8).Range.Text = "0.204"
$t.Cell(8, 4).Range.Text = "29 (38.7)"
$t.Cell(9, 4).Range.Text = "51 (68.0)"
$t.Cell(9, 8).Range.Text = "0.708"
$t.Cell(10, 4).Range.Text = "24 (32.0)"
$t.Cell(11, 4).Range.Text = "75 (100.0)"
$t.Cell(12, 4).Range.Text = "75 (100.0)"
$t.Cell(13, 4).Range.Text = "75 (100.0)"
$t.Cell(13, 6).Range.Text = "0.725"
$t.Cell(15, 4).Range.Text = "73 (97.3)"
$t.Cell(15, 6).Range.Text = "0.251"
$t.Cell(15, 8).Range.Text = "0.437"
$t.Cell(16, 4).Range.Text = "2 (2.7)"
$t.Cell(17, 4).Range.Text = "73 (97.3)"
$t.Cell(17, 6).Range.Text = "0.581"
$t.Cell(17, 8).Range.Text = "0.986"
$t.Cell(18, 4).Range.Text = "2 (2.7)"
$t.Cell(19, 4).Range.Text = "70 (93.3)"
$t.Cell(19, 6).Range.Text = "0.137"
$t.Cell(19, 8).Range.Text = "0.278"
$t.Cell(20, 4).Range.Text = "5 (6.7)"
$t.Cell(21, 4).Range.Text = "75 (100.0)"
$t.Cell(22, 4).Range.Text = "74 (98.7)"
$t.Cell(22, 6).Range.Text = "0.714"
$t.Cell(24, 4).Range.Text = "74 (98.7)"
$t.Cell(24, 6).Range.Text = "0.680"
$t.Cell(26, 4).Range.Text = "71 (94.7)"
$t.Cell(26, 6).Range.Text = "0.320"
$t.Cell(28, 4).Range.Text = "71 (94.7)"
$t.Cell(28, 6).Range.Text = "0.359"
$t.Cell(28, 8).Range.Text = "0.872"
$t.Cell(30, 4).Range.Text = "75 (100.0)"
$t.Cell(30, 6).Range.Text = "0.379"
$t.Cell(30, 8).Range.Text = "0.621"
$t.Cell(32, 4).Range.Text = "75 (100.0)"
$t.Cell(32, 6).Range.Text = "0.725"
$t.Cell(34, 4).Range.Text = "75 (100.0)"
$t.Cell(35, 4).Range.Text = "74 (98.7)"
$t.Cell(35, 6).Range.Text = "0.783"
$t.Cell(37, 4).Range.Text = "75 (100.0)"
$t.Cell(37, 6).Range.Text = "0.725"
$t.Cell(39, 4).Range.Text = "75 (100.0)"
$t.Cell(40, 4).Range.Text = "75 (100.0)"
$t.Cell(41, 4).Range.Text = "73 (97.3)"
$t.Cell(41, 6).Range.Text = "0.251"
$t.Cell(41, 8).Range.Text = "0.437"
$t.Cell(42, 4).Range.Text = "2 (2.7)"
$t.Cell(43, 4).Range.Text = "75 (100.0)"
$t.Cell(43, 6).Range.Text = "0.273"
$t.Cell(43, 8).Range.Text = "0.449"
$t.Cell(45, 4).Range.Text = "37 (49.3)"
$t.Cell(45, 8).Range.Text = "0.454"
$t.Cell(46, 4).Range.Text = "38 (50.7)"
$t.Cell(47, 4).Range.Text = "37 (49.3)"
$t.Cell(47, 8).Range.Text = "0.477"
$t.Cell(48, 4).Range.Text = "16 (21.3)"
$t.Cell(50, 4).Range.Text = "12 (16.0)"
$t.Cell(51, 4).Range.Text = "10 (13.3)"
$t.Cell(52, 4).Range.Text = "24 (32.0)"
$t.Cell(52, 8).Range.Text = "0.452"
$t.Cell(53, 4).Range.Text = "51 (68.0)"
$t.Cell(54, 4).Range.Text = "24 (32.0)"
$t.Cell(54, 8).Range.Text = "0.486"
$t.Cell(55, 4).Range.Text = "21 (28.0)"
$t.Cell(56, 4).Range.Text = "3 (4.0)"
$t.Cell(57, 4).Range.Text = "18 (24.0)"
$t.Cell(58, 4).Range.Text = "9 (12.0)"
$t.Cell(59, 4).Range.Text = "70 (93.3)"
$t.Cell(60, 4).Range.Text = "5 (6.7)"
$t.Cell(61, 4).Range.Text = "70 (93.3)"
$t.Cell(61, 6).Range.Text = "0.588"
$t.Cell(61, 8).Range.Text = "0.542"
$t.Cell(66, 4).Range.Text = "74 (98.7)"
$t.Cell(66, 6).Range.Text = "0.714"
$t.Cell(68, 4).Range.Text = "73 (97.3)"
$t.Cell(68, 6).Range.Text = "0.856"
$t.Cell(68, 8).Range.Text = "0.986"
$t.Cell(69, 4).Range.Text = "2 (2.7)"
$t.Cell(70, 4).Range.Text = "63 (84.0)"
$t.Cell(70, 6).Range.Text = "0.273"
$t.Cell(70, 8).Range.Text = "0.239"
$t.Cell(71, 4).Range.Text = "12 (16.0)"
$t.Cell(72, 4).Range.Text = "75 (100.0)"
$t.Cell(72, 6).Range.Text = "0.273"
$t.Cell(72, 8).Range.Text = "0.449"
